$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1) - F column "想去人数" updates
$wsExhibit = $wb.Worksheets.Item(1)
$wsExhibit.Range("F3").Value = 3814
$wsExhibit.Range("F4").Value = 2284
$wsExhibit.Range("F5").Value = 451
$wsExhibit.Range("F7").Value = 20
$wsExhibit.Range("F8").Value = 185
$wsExhibit.Range("F9").Value = 109
$wsExhibit.Range("F10").Value = 98
$wsExhibit.Range("F11").Value = 1417
$wsExhibit.Range("F12").Value = 248
$wsExhibit.Range("F13").Value = 2410
$wsExhibit.Range("F14").Value = 171

# Sheet "全部类型" (index 4) - F column "想去人数" updates
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F3").Value = 3814
$wsAll.Range("F4").Value = 2284
$wsAll.Range("F5").Value = 451
$wsAll.Range("F7").Value = 20
$wsAll.Range("F9").Value = 185
$wsAll.Range("F10").Value = 109
$wsAll.Range("F11").Value = 98
$wsAll.Range("F14").Value = 1417
$wsAll.Range("F15").Value = 248
$wsAll.Range("F16").Value = 2410
$wsAll.Range("F17").Value = 171
